$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 688.4583
$ws.Range("I53").Value = 223.21428
$ws.Range("J53").Value = 1339.8
$ws.Range("K53").Value = 223.21428
$ws.Range("L53").Value = 1339.8
$ws.Range("M53").Value = 413.78572
$ws.Range("N53").Value = -2613.8
$ws.Range("H55").Value = 203.47058
$ws.Range("J55").Value = 241.81818
$ws.Range("L55").Value = 241.81818
$ws.Range("N55").Value = -669.81818
$ws.Range("H112").Value = 6026252.5
$ws.Range("J112").Value = 6331238.5
$ws.Range("L112").Value = 18993715.5
$ws.Range("N112").Value = -18995931.5
$ws.Range("H129").Value = 2861.5
$ws.Range("I129").Value = 1117.4
$ws.Range("J129").Value = 4605.6
$ws.Range("K129").Value = 3352.2
$ws.Range("L129").Value = 13816.8
$ws.Range("M129").Value = 1647.8
$ws.Range("N129").Value = -23816.8
$ws.Range("H137").Value = 42555120
$ws.Range("I137").Value = 25001962
$ws.Range("J137").Value = 142858880
$ws.Range("K137").Value = 75005886
$ws.Range("L137").Value = 428576640
$ws.Range("M137").Value = -75003336
$ws.Range("N137").Value = -428581740
$ws.Range("H138").Value = 7579048.5
$ws.Range("I138").Value = 2834.75
$ws.Range("J138").Value = 8336670
$ws.Range("K138").Value = 8504.25
$ws.Range("L138").Value = 25010010
$ws.Range("M138").Value = -3364.25
$ws.Range("N138").Value = -25020290
$ws.Range("H141").Value = 2388.9375
$ws.Range("I141").Value = 2494.6667
$ws.Range("K141").Value = 7484.000100000001
$ws.Range("M141").Value = -2304.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2283.5908
$ws.Range("I2").Value = 2416.2666
$ws.Range("K2").Value = 2416.2666
$ws.Range("M2").Value = -2303.2666
$ws.Range("H32").Value = 10424104
$ws.Range("I32").Value = 12990309
$ws.Range("J32").Value = 24220.21
$ws.Range("K32").Value = 12990309
$ws.Range("L32").Value = 24220.21
$ws.Range("M32").Value = -12990022
$ws.Range("N32").Value = -24794.21
$ws.Range("H61").Value = 26320800
$ws.Range("I61").Value = 40003430
$ws.Range("K61").Value = 40003430
$ws.Range("M61").Value = -40003218
$ws.Range("H116").Value = 2283.5908
$ws.Range("I116").Value = 2416.2666
$ws.Range("K116").Value = 2416.2666
$ws.Range("M116").Value = -122.2665999999999
$ws.Range("H136").Value = 26320800
$ws.Range("I136").Value = 40003430
$ws.Range("K136").Value = 120010290
$ws.Range("M136").Value = -120007740

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2283.5908
$ws.Range("I3").Value = 2416.2666
$ws.Range("K3").Value = 2416.2666
$ws.Range("M3").Value = -2302.2666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 102013300
$ws.Range("J4").Value = 200000400
$ws.Range("L4").Value = 200000400
$ws.Range("N4").Value = -200000624
$ws.Range("H58").Value = 1317
$ws.Range("I58").Value = 984.7646999999999
$ws.Range("J58").Value = 3199.6667
$ws.Range("K58").Value = 984.7646999999999
$ws.Range("L58").Value = 3199.6667
$ws.Range("M58").Value = -781.7646999999999
$ws.Range("N58").Value = -3605.6667
$ws.Range("H68").Value = 39765.832
$ws.Range("J68").Value = 39765.832
$ws.Range("L68").Value = 39765.832
$ws.Range("N68").Value = -41263.832
$ws.Range("H71").Value = 39765.832
$ws.Range("J71").Value = 39765.832
$ws.Range("L71").Value = 119297.496
$ws.Range("N71").Value = -126785.496
$ws.Range("H74").Value = 39285
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 39285
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H86").Value = 3685.6428
$ws.Range("I86").Value = 2838.125
$ws.Range("K86").Value = 2838.125
$ws.Range("M86").Value = -1715.125
$ws.Range("H89").Value = 3685.6428
$ws.Range("I89").Value = 2838.125
$ws.Range("K89").Value = 14190.625
$ws.Range("M89").Value = -8574.625
$ws.Range("H132").Value = 4287.75
$ws.Range("I132").Value = 3491.6667
$ws.Range("J132").Value = 6676
$ws.Range("K132").Value = 10475.0001
$ws.Range("L132").Value = 20028
$ws.Range("M132").Value = -7945.000100000001
$ws.Range("N132").Value = -25088
$ws.Range("H136").Value = 1317
$ws.Range("I136").Value = 984.7646999999999
$ws.Range("J136").Value = 3199.6667
$ws.Range("K136").Value = 2954.2941
$ws.Range("L136").Value = 9599.000100000001
$ws.Range("M136").Value = -404.2941000000001
$ws.Range("N136").Value = -14699.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34965572
$ws.Range("I4").Value = 66435124
$ws.Range("K4").Value = 199305372
$ws.Range("M4").Value = -199305260
$ws.Range("H113").Value = 3825.7856
$ws.Range("I113").Value = 3749.5
$ws.Range("J113").Value = 3838.5
$ws.Range("K113").Value = 11248.5
$ws.Range("L113").Value = 11515.5
$ws.Range("M113").Value = -9078.5
$ws.Range("N113").Value = -15855.5
$ws.Range("H131").Value = 40375.62
$ws.Range("I131").Value = 81045.38
$ws.Range("J131").Value = 7331.4375
$ws.Range("K131").Value = 243136.14
$ws.Range("L131").Value = 21994.3125
$ws.Range("M131").Value = -238096.14
$ws.Range("N131").Value = -32074.3125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 287.4643
$ws.Range("I2").Value = 65.9375
$ws.Range("J2").Value = 582.8333
$ws.Range("K2").Value = 65.9375
$ws.Range("L2").Value = 582.8333
$ws.Range("M2").Value = 47.0625
$ws.Range("N2").Value = -808.8333
$ws.Range("H20").Value = 23492.5
$ws.Range("J20").Value = 23492.5
$ws.Range("L20").Value = 23492.5
$ws.Range("N20").Value = -23982.5
$ws.Range("H24").Value = 27505.21
$ws.Range("I24").Value = 14999.667
$ws.Range("K24").Value = 14999.667
$ws.Range("M24").Value = -14826.667
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H70").Value = 5802.3335
$ws.Range("I70").Value = 4119.2856
$ws.Range("K70").Value = 4119.2856
$ws.Range("M70").Value = -3849.2856
$ws.Range("H73").Value = 5802.3335
$ws.Range("I73").Value = 4119.2856
$ws.Range("K73").Value = 4119.2856
$ws.Range("M73").Value = -3183.2856
$ws.Range("H107").Value = 936.63635
$ws.Range("I107").Value = 1169.4
$ws.Range("J107").Value = 742.6667
$ws.Range("K107").Value = 1169.4
$ws.Range("L107").Value = 742.6667
$ws.Range("M107").Value = 750.5999999999999
$ws.Range("N107").Value = -4582.6667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 511.1154
$ws.Range("I55").Value = 453.42105
$ws.Range("K55").Value = 453.42105
$ws.Range("M55").Value = -280.42105
$ws.Range("H100").Value = 2814.0625
$ws.Range("I100").Value = 1840.5
$ws.Range("J100").Value = 3787.625
$ws.Range("K100").Value = 1840.5
$ws.Range("L100").Value = 3787.625
$ws.Range("M100").Value = -1299.5
$ws.Range("N100").Value = -4869.625
$ws.Range("H132").Value = 68966760
$ws.Range("I132").Value = 1242.5
$ws.Range("J132").Value = 222223470
$ws.Range("K132").Value = 3727.5
$ws.Range("L132").Value = 666670410
$ws.Range("M132").Value = -1197.5
$ws.Range("N132").Value = -666675470
$ws.Range("H136").Value = 3444.1667
$ws.Range("I136").Value = 3444.1667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10332.5001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7782.500100000001
$ws.Range("N136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1075
$ws.Range("J8").Value = 1075
$ws.Range("L8").Value = 1075
$ws.Range("N8").Value = -1355
$ws.Range("H38").Value = 20000000
$ws.Range("I38").Value = 20000000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 20000000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -19999527
$ws.Range("N38").ClearContents()
$ws.Range("H136").Value = 1991.1395
$ws.Range("I136").Value = 1868.7805
$ws.Range("J136").Value = 4499.5
$ws.Range("K136").Value = 5606.3415
$ws.Range("L136").Value = 13498.5
$ws.Range("M136").Value = -3056.3415
$ws.Range("N136").Value = -18598.5
